# Actualizacion automatica: inserta una nueva fila para la clienta
# "LINCANGO LUGMANIA SANDY LIZETH" justo antes de la fila de
# "VACA PANCHI DORYS CAROLINA" (antigua fila 5) en ambas hojas, y
# actualiza las filas de totales ("0 de N") para reflejar el nuevo
# numero de clientes.

$wb = $excel.ActiveWorkbook

# --- Hoja 1: VENTAS POR GRUPO (columnas A:N) ---
$ws1 = $wb.Worksheets.Item(1)

# Inserta una fila nueva en la posicion 5, desplazando hacia abajo la
# fila existente (VACA PANCHI DORYS CAROLINA) y la fila de totales.
$ws1.Rows.Item(5).Insert()

$ws1.Cells.Item(5, 1).Value = "VACA PANCHI CAROLINA"
$ws1.Cells.Item(5, 2).Value = "LINCANGO LUGMANIA SANDY LIZETH"
for ($col = 3; $col -le 14; $col++) {
    $ws1.Cells.Item(5, $col).Value = 0
}

# La fila de totales (ahora en la fila 7) pasa de "0 de 4" a "0 de 5".
for ($col = 3; $col -le 14; $col++) {
    $ws1.Cells.Item(7, $col).Value = "0 de 5"
}

# --- Hoja 2: VENTA MENSUAL (columnas A:F) ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(5).Insert()

$ws2.Cells.Item(5, 1).Value = "VACA PANCHI CAROLINA"
$ws2.Cells.Item(5, 2).Value = "LINCANGO LUGMANIA SANDY LIZETH"
for ($col = 3; $col -le 6; $col++) {
    $ws2.Cells.Item(5, $col).Value = 0
}
